$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2, F2, G2
$ws.Range("D2").Value = 26044
$ws.Range("E2").Value = 506958697449
$ws.Range("F2").Value = 3434370591
$ws.Range("G2").Value = 0.07351000000000001

# Row 3: update D3, E3, F3, G3
$ws.Range("D3").Value = 1648.51
$ws.Range("E3").Value = 198127113636
$ws.Range("F3").Value = 2370169187
$ws.Range("G3").Value = -0.02006

# Row 4: update D4, E4, F4, G4
$ws.Range("D4").Value = 0.999378
$ws.Range("E4").Value = 82802491526
$ws.Range("F4").Value = 5828650468
$ws.Range("G4").Value = 0.01024

# Row 5: update D5, E5, F5, G5
$ws.Range("D5").Value = 216.72
$ws.Range("E5").Value = 33329943321
$ws.Range("F5").Value = 190237548
$ws.Range("G5").Value = -0.06942

# Row 6: update D6, E6, F6, G6
$ws.Range("D6").Value = 0.526204
$ws.Range("E6").Value = 27847136797
$ws.Range("F6").Value = 466732452
$ws.Range("G6").Value = 0.77638

# Row 7: update D7, E7, F7, G7
$ws.Range("D7").Value = 0.999864
$ws.Range("E7").Value = 26018604483
$ws.Range("F7").Value = 1445591984
$ws.Range("G7").Value = -0.04011

# Row 8: update D8, E8, F8, G8
$ws.Range("D8").Value = 1647.32
$ws.Range("E8").Value = 13834721080
$ws.Range("F8").Value = 3947369
$ws.Range("G8").Value = -0.10148

# Row 9: update D9, E9, F9, G9
$ws.Range("D9").Value = 0.263356
$ws.Range("E9").Value = 9224372100
$ws.Range("F9").Value = 127404453
$ws.Range("G9").Value = 1.74286

# Row 10: update D10, E10, F10, G10
$ws.Range("D10").Value = 0.063135
$ws.Range("E10").Value = 8882092728
$ws.Range("F10").Value = 140667452
$ws.Range("G10").Value = 0.54836

# Row 11: update D11, E11, F11, G11
$ws.Range("D11").Value = 20.6
$ws.Range("E11").Value = 8389845541
$ws.Range("F11").Value = 163396346
$ws.Range("G11").Value = 1.79039

# Row 12: update D12, E12, F12, G12
$ws.Range("D12").Value = 0.077309
$ws.Range("E12").Value = 6908836430
$ws.Range("F12").Value = 129835797
$ws.Range("G12").Value = -0.15625

# Row 13: update D13, E13, F13, G13
$ws.Range("D13").Value = 4.51
$ws.Range("E13").Value = 5718452130
$ws.Range("F13").Value = 60712450
$ws.Range("G13").Value = 0.1519

# Row 14: update B14, C14, D14, E14, F14, G14
$ws.Range("B14").Value = "TON"
$ws.Range("C14").Value = "Toncoin"
$ws.Range("D14").Value = 1.48
$ws.Range("E14").Value = 5102228069
$ws.Range("F14").Value = 22210101
$ws.Range("G14").Value = 3.64884

# Row 15: update B15, C15, D15, E15, F15, G15
$ws.Range("B15").Value = "MATIC"
$ws.Range("C15").Value = "Polygon"
$ws.Range("D15").Value = 0.547096
$ws.Range("E15").Value = 5096431700
$ws.Range("F15").Value = 130975474
$ws.Range("G15").Value = 0.01687

# Row 16: update B16, C16, D16, E16, F16, G16
$ws.Range("B16").Value = "LTC"
$ws.Range("C16").Value = "Litecoin"
$ws.Range("D16").Value = 65.2
$ws.Range("E16").Value = 4796859420
$ws.Range("F16").Value = 188968239
$ws.Range("G16").Value = -0.28737

# Row 17: update B17, C17, D17, E17, F17, G17
$ws.Range("B17").Value = "SHIB"
$ws.Range("C17").Value = "Shiba Inu"
$ws.Range("D17").Value = 0.00000808
$ws.Range("E17").Value = 4767444805
$ws.Range("F17").Value = 73137453
$ws.Range("G17").Value = -1.15629

# Row 18: update D18, E18, F18, G18
$ws.Range("D18").Value = 26049
$ws.Range("E18").Value = 4240351254
$ws.Range("F18").Value = 23139057
$ws.Range("G18").Value = 0.1081

# Row 19: update D19, E19, F19, G19
$ws.Range("D19").Value = 0.999682
$ws.Range("E19").Value = 3895962922
$ws.Range("F19").Value = 47510977
$ws.Range("G19").Value = -0.03604

# Row 20: update B20, C20, D20, E20, F20, G20
$ws.Range("B20").Value = "BCH"
$ws.Range("C20").Value = "Bitcoin Cash"
$ws.Range("D20").Value = 191.16
$ws.Range("E20").Value = 3723021445
$ws.Range("F20").Value = 68082496
$ws.Range("G20").Value = 0.93882

# Row 21: update B21, C21, D21, E21, F21, G21
$ws.Range("B21").Value = "AVAX"
$ws.Range("C21").Value = "Avalanche"
$ws.Range("D21").Value = 10.18
$ws.Range("E21").Value = 3598828249
$ws.Range("F21").Value = 92423171
$ws.Range("G21").Value = 1.40042

# Row 22: update B22, C22, D22, E22, F22, G22
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "LEO Token"
$ws.Range("D22").Value = 3.81
$ws.Range("E22").Value = 3540009135
$ws.Range("F22").Value = 142091
$ws.Range("G22").Value = -1.19456

# Row 23: update B23, C23, D23, E23, F23, G23
$ws.Range("B23").Value = "UNI"
$ws.Range("C23").Value = "Uniswap"
$ws.Range("D23").Value = 4.63
$ws.Range("E23").Value = 3486187353
$ws.Range("F23").Value = 56023469
$ws.Range("G23").Value = 1.53667

# Row 24: update D24, E24, F24, G24
$ws.Range("D24").Value = 0.121001
$ws.Range("E24").Value = 3316709829
$ws.Range("F24").Value = 55830039
$ws.Range("G24").Value = -1.53255

# Row 25: update D25, E25, F25, G25
$ws.Range("D25").Value = 5.97
$ws.Range("E25").Value = 3207520741
$ws.Range("F25").Value = 88308038
$ws.Range("G25").Value = -0.648

# Row 26: update D26, E26, F26, G26
$ws.Range("D26").Value = 0.99984
$ws.Range("E26").Value = 3144617444
$ws.Range("F26").Value = 601426693
$ws.Range("G26").Value = 0.00784

# Row 27: update D27, E27, F27, G27
$ws.Range("D27").Value = 0.998996
$ws.Range("E27").Value = 2904170794
$ws.Range("F27").Value = 470443550
$ws.Range("G27").Value = 0.03631

# Row 28: update D28, E28, F28, G28
$ws.Range("D28").Value = 144.34
$ws.Range("E28").Value = 2619510023
$ws.Range("F28").Value = 70241292
$ws.Range("G28").Value = 1.2331

# Row 29: update D29, E29, F29, G29
$ws.Range("D29").Value = 42.87
$ws.Range("E29").Value = 2572151879
$ws.Range("F29").Value = 1535814
$ws.Range("G29").Value = 0.03573

# Row 30: update B30, C30, D30, E30, F30, G30
$ws.Range("B30").Value = "ETC"
$ws.Range("C30").Value = "Ethereum Classic"
$ws.Range("D30").Value = 15.85
$ws.Range("E30").Value = 2260260690
$ws.Range("F30").Value = 41445753
$ws.Range("G30").Value = 0.1097

# Row 31: update B31, C31, D31, E31, F31, G31
$ws.Range("B31").Value = "ATOM"
$ws.Range("C31").Value = "Cosmos Hub"
$ws.Range("D31").Value = 7.12
$ws.Range("E31").Value = 2081985459
$ws.Range("F31").Value = 77102043
$ws.Range("G31").Value = -1.18019

# Row 32: update D32, E32, F32, G32
$ws.Range("D32").Value = 0.057744
$ws.Range("E32").Value = 1907007600
$ws.Range("F32").Value = 54014405
$ws.Range("G32").Value = 0.06707

# Row 33: update B33, C33, D33, E33, F33, G33
$ws.Range("B33").Value = "ICP"
$ws.Range("C33").Value = "Internet Computer"
$ws.Range("D33").Value = 3.54
$ws.Range("E33").Value = 1562823720
$ws.Range("F33").Value = 16920048
$ws.Range("G33").Value = -0.0267

# Row 34: update B34, C34, D34, E34, F34, G34
$ws.Range("B34").Value = "QNT"
$ws.Range("C34").Value = "Quant"
$ws.Range("D34").Value = 103.31
$ws.Range("E34").Value = 1504118148
$ws.Range("F34").Value = 15589299
$ws.Range("G34").Value = -1.24585

# Row 35: update B35, C35, D35, E35, F35, G35
$ws.Range("B35").Value = "FIL"
$ws.Range("C35").Value = "Filecoin"
$ws.Range("D35").Value = 3.24
$ws.Range("E35").Value = 1436290899
$ws.Range("F35").Value = 53987319
$ws.Range("G35").Value = -0.50931

# Row 36: update B36, C36, D36, E36, F36, G36
$ws.Range("B36").Value = "MNT"
$ws.Range("C36").Value = "Mantle"
$ws.Range("D36").Value = 0.435081
$ws.Range("E36").Value = 1407806378
$ws.Range("F36").Value = 3310735
$ws.Range("G36").Value = 0.59217

# Row 37: update B37, C37, D37, E37, F37, G37
$ws.Range("B37").Value = "LDO"
$ws.Range("C37").Value = "Lido DAO"
$ws.Range("D37").Value = 1.59
$ws.Range("E37").Value = 1405900344
$ws.Range("F37").Value = 45896445
$ws.Range("G37").Value = -0.05311

# Row 38: update D38, E38, F38, G38
$ws.Range("D38").Value = 0.05107
$ws.Range("E38").Value = 1339413478
$ws.Range("F38").Value = 3029891
$ws.Range("G38").Value = 0.00362

# Row 39: update D39, E39, F39, G39
$ws.Range("D39").Value = 5.7
$ws.Range("E39").Value = 1297184221
$ws.Range("F39").Value = 29566935
$ws.Range("G39").Value = -0.79947

# Row 40: update D40, E40, F40, G40
$ws.Range("D40").Value = 0.956992
$ws.Range("E40").Value = 1219301114
$ws.Range("F40").Value = 90823071
$ws.Range("G40").Value = 1.82277

# Row 41: update D41, E41, F41, G41
$ws.Range("D41").Value = 0.01584905
$ws.Range("E41").Value = 1151043225
$ws.Range("F41").Value = 19879938
$ws.Range("G41").Value = -0.77488

# Row 42: update D42, E42, F42, G42
$ws.Range("D42").Value = 1.2
$ws.Range("E42").Value = 1125050969
$ws.Range("F42").Value = 53752245
$ws.Range("G42").Value = -0.50285

# Row 43: update E43, F43, G43
$ws.Range("E43").Value = 1051579744
$ws.Range("F43").Value = 58995108
$ws.Range("G43").Value = -1.51614

# Row 44: update D44, E44, F44, G44
$ws.Range("D44").Value = 1053.26
$ws.Range("E44").Value = 948356738
$ws.Range("F44").Value = 42789829
$ws.Range("G44").Value = 2.21602

# Row 45: update D45, E45, F45, G45
$ws.Range("D45").Value = 1785.7
$ws.Range("E45").Value = 917437844
$ws.Range("F45").Value = 11294099
$ws.Range("G45").Value = -0.24803

# Row 46: update B46, C46, D46, E46, F46, G46
$ws.Range("B46").Value = "AAVE"
$ws.Range("C46").Value = "Aave"
$ws.Range("D46").Value = 57.12
$ws.Range("E46").Value = 828438467
$ws.Range("F46").Value = 51323287
$ws.Range("G46").Value = 0.35197

# Row 47: update B47, C47, D47, E47, F47, G47
$ws.Range("B47").Value = "GRT"
$ws.Range("C47").Value = "The Graph"
$ws.Range("D47").Value = 0.08867
$ws.Range("E47").Value = 811983493
$ws.Range("F47").Value = 25680687
$ws.Range("G47").Value = -0.79378

# Row 48: update B48, C48, D48, E48, F48, G48
$ws.Range("B48").Value = "FRAX"
$ws.Range("C48").Value = "Frax"
$ws.Range("D48").Value = 0.997003
$ws.Range("E48").Value = 804268385
$ws.Range("F48").Value = 3680667
$ws.Range("G48").Value = -0.11491

# Row 49: update B49, C49, D49, E49, F49, G49
$ws.Range("B49").Value = "WBT"
$ws.Range("C49").Value = "WhiteBIT Coin"
$ws.Range("D49").Value = 5.37
$ws.Range("E49").Value = 774075836
$ws.Range("F49").Value = 6518263
$ws.Range("G49").Value = -0.27053

# Row 50: update B50, C50, D50, E50, F50, G50
$ws.Range("B50").Value = "ALGO"
$ws.Range("C50").Value = "Algorand"
$ws.Range("D50").Value = 0.094763
$ws.Range("E50").Value = 741014084
$ws.Range("F50").Value = 13647532
$ws.Range("G50").Value = -1.01591

# Row 51: update B51, C51, D51, E51, F51, G51
$ws.Range("B51").Value = "KAS"
$ws.Range("C51").Value = "Kaspa"
$ws.Range("D51").Value = 0.03670698
$ws.Range("E51").Value = 737570345
$ws.Range("F51").Value = 8085140
$ws.Range("G51").Value = 2.14587
